$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- old row 9
$ws.Range("D2").Value2 = 44630
$ws.Range("J2").Value2 = 90
$ws.Range("K2").Value2 = 2500
$ws.Range("L2").Value2 = 3000
$ws.Range("M2").Value2 = 2722
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value2 = 454

# Row 4 <- old row 11
$ws.Range("D4").Value2 = 44987
$ws.Range("J4").Value2 = 130
$ws.Range("K4").Value2 = 4500
$ws.Range("L4").Value2 = 5000
$ws.Range("M4").Value2 = 4692
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value2 = 782

# Row 5 <- old row 7
$ws.Range("D5").Value2 = 44643
$ws.Range("J5").Value2 = 90
$ws.Range("K5").Value2 = 2800
$ws.Range("L5").Value2 = 3000
$ws.Range("M5").Value2 = 2911
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value2 = 485

# Row 6 <- old row 14
$ws.Range("D6").Value2 = 44658
$ws.Range("J6").Value2 = 180
$ws.Range("K6").Value2 = 2500
$ws.Range("L6").Value2 = 3000
$ws.Range("M6").Value2 = 2778
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value2 = 463

# Row 7 <- old row 4
$ws.Range("D7").Value2 = 44876
$ws.Range("J7").Value2 = 80
$ws.Range("K7").Value2 = 6500
$ws.Range("L7").Value2 = 7000
$ws.Range("M7").Value2 = 6812
$ws.Range("O7").Value = "Región Metropolitana"
$ws.Range("P7").Value2 = 1135

# Row 8 <- old row 15
$ws.Range("D8").Value2 = 44631
$ws.Range("J8").Value2 = 110
$ws.Range("K8").Value2 = 3000
$ws.Range("L8").Value2 = 3500
$ws.Range("M8").Value2 = 3273
$ws.Range("O8").Value = "Provincia de Chacabuco"
$ws.Range("P8").Value2 = 546

# Row 9 <- old row 13
$ws.Range("D9").Value2 = 44685
$ws.Range("J9").Value2 = 150
$ws.Range("K9").Value2 = 3000
$ws.Range("L9").Value2 = 3500
$ws.Range("M9").Value2 = 3267
$ws.Range("O9").Value = "Región Metropolitana"
$ws.Range("P9").Value2 = 544

# Row 10 <- old row 8
$ws.Range("D10").Value2 = 45021
$ws.Range("J10").Value2 = 50
$ws.Range("K10").Value2 = 4500
$ws.Range("L10").Value2 = 5000
$ws.Range("M10").Value2 = 4700
$ws.Range("O10").Value = "Región Metropolitana"
$ws.Range("P10").Value2 = 783

# Row 11 <- old row 16
$ws.Range("D11").Value2 = 44644
$ws.Range("J11").Value2 = 140
$ws.Range("K11").Value2 = 2500
$ws.Range("L11").Value2 = 3000
$ws.Range("M11").Value2 = 2786
$ws.Range("O11").Value = "Provincia de Chacabuco"
$ws.Range("P11").Value2 = 464

# Row 12 <- old row 5
$ws.Range("D12").Value2 = 44957
$ws.Range("J12").Value2 = 70
$ws.Range("K12").Value2 = 1500
$ws.Range("L12").Value2 = 2000
$ws.Range("M12").Value2 = 1857
$ws.Range("O12").Value = "Región Metropolitana"
$ws.Range("P12").Value2 = 310

# Row 13 <- old row 2
$ws.Range("D13").Value2 = 44672
$ws.Range("J13").Value2 = 140
$ws.Range("K13").Value2 = 3000
$ws.Range("L13").Value2 = 3500
$ws.Range("M13").Value2 = 3286
$ws.Range("O13").Value = "Región Metropolitana"
$ws.Range("P13").Value2 = 548

# Row 14 <- old row 10
$ws.Range("D14").Value2 = 44650
$ws.Range("J14").Value2 = 130
$ws.Range("K14").Value2 = 3000
$ws.Range("L14").Value2 = 3500
$ws.Range("M14").Value2 = 3308
$ws.Range("O14").Value = "Región Metropolitana"
$ws.Range("P14").Value2 = 551

# Row 15 <- old row 12
$ws.Range("D15").Value2 = 44659
$ws.Range("J15").Value2 = 90
$ws.Range("K15").Value2 = 2500
$ws.Range("L15").Value2 = 3000
$ws.Range("M15").Value2 = 2722
$ws.Range("O15").Value = "Región Metropolitana"
$ws.Range("P15").Value2 = 454

# Row 16 <- old row 6
$ws.Range("D16").Value2 = 44637
$ws.Range("J16").Value2 = 170
$ws.Range("K16").Value2 = 2800
$ws.Range("L16").Value2 = 3000
$ws.Range("M16").Value2 = 2906
$ws.Range("O16").Value = "Región Metropolitana"
$ws.Range("P16").Value2 = 484
